$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 338, pushing the existing data (rows 338-358) down to
# rows 339-359, and grow the used range to A1:T359.
$ws.Rows.Item(338).Insert()

# Populate the newly inserted row 338 with the new weekly record.
$ws.Range("A338").Value = 6
$ws.Range("B338").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C338").Value = "Metropolitana"
$ws.Range("D338").Value = 45041
$ws.Range("E338").Value = 13
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100101
$ws.Range("H338").Value = "Berries"
$ws.Range("I338").Value = 100101004
$ws.Range("J338").Value = "Frambuesa"
$ws.Range("K338").Value = "Sin especificar"
$ws.Range("L338").Value = "Primera"
$ws.Range("M338").Value = 75
$ws.Range("N338").Value = 9000
$ws.Range("O338").Value = 9000
$ws.Range("P338").Value = 9000
$ws.Range("Q338").Value = '$/bandeja 2 kilos'
$ws.Range("R338").Value = "Provincia de Curicó"
$ws.Range("S338").Value = 4500
$ws.Range("T338").Value = 2
